# Refresh the cryptos table (Coin/Link/Price/Volume) to the latest snapshot.
# Source: scheduled GitHub Actions data pull (cryptos.xlsx).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "44.617.98"
$ws.Range("E2").Value = "  +3.77%  "
# Row 3
$ws.Range("D3").Value = "2.418.23"
$ws.Range("E3").Value = "  +2.10%  "
# Row 4
$ws.Range("E4").Value = "  +0.01%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.78"
$ws.Range("E5").Value = "  +4.25%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.19"
$ws.Range("E6").Value = "  +6.03%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.512"
$ws.Range("E7").Value = "  +2.01%  "
# Row 8
$ws.Range("E8").Value = "  -0.04%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("E9").Value = "  +10.36%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.26"
$ws.Range("E10").Value = "  +2.65%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0798"
$ws.Range("E11").Value = "  +1.52%  "
# Row 12
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.122"
$ws.Range("E12").Value = "  -1.62%  "
# Row 13
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.60"
$ws.Range("E13").Value = "  +1.42%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.89"
$ws.Range("E14").Value = "  +2.50%  "
# Row 15
$ws.Range("D15").Value = "2.797.39"
$ws.Range("E15").Value = "  +2.40%  "
# Row 16
$ws.Range("D16").Value = "2.429.14"
$ws.Range("E16").Value = "  +3.06%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.828"
$ws.Range("E17").Value = "  +4.26%  "
# Row 18
$ws.Range("D18").Value = "44.485.03"
$ws.Range("E18").Value = "  +3.55%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.18"
$ws.Range("E19").Value = "  +1.76%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.34"
$ws.Range("E20").Value = "  +1.25%  "
# Row 21
$ws.Range("D21").Value = "0.0₃0916"
$ws.Range("E21").Value = "  +3.46%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.59"
$ws.Range("E22").Value = "  +1.01%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.74"
$ws.Range("E23").Value = "  +2.75%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.27"
$ws.Range("E24").Value = "  +3.52%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.49"
$ws.Range("E25").Value = "  +2.23%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.06%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.18"
$ws.Range("E27").Value = "  +3.09%  "
# Row 28
$ws.Range("E28").Value = "  -3.43%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.49"
$ws.Range("E29").Value = "  +2.16%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.27"
$ws.Range("E30").Value = "  +2.92%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "48.32"
$ws.Range("E31").Value = "  +1.16%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.126"
$ws.Range("E32").Value = "  +16.62%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.49"
$ws.Range("E33").Value = "  +10.97%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.15"
$ws.Range("E34").Value = "  +2.84%  "
# Row 35
$ws.Range("E35").Value = "  +0.28%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0763"
$ws.Range("E36").Value = "  +4.95%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.88"
$ws.Range("E37").Value = "  +2.77%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.45"
$ws.Range("E38").Value = "  +2.85%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.84"
$ws.Range("E39").Value = "  +0.14%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "124.71"
$ws.Range("E40").Value = "  -3.45%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.109"
$ws.Range("E41").Value = "  +1.59%  "
# Row 42
$ws.Range("E42").Value = "  -3.44%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.83"
$ws.Range("E43").Value = "  +0.13%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0289"
$ws.Range("E44").Value = "  +3.67%  "
# Row 45
$ws.Range("D45").Value = "1.939.86"
$ws.Range("E45").Value = "  +0.70%  "
# Row 46
$ws.Range("E46").Value = "  -1.38%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.92"
$ws.Range("E47").Value = "  +7.76%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.16"
$ws.Range("E48").Value = "  -0.64%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.73"
$ws.Range("E49").Value = "  +14.99%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.02"
$ws.Range("E50").Value = "  +5.07%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.78"
$ws.Range("E51").Value = "  +5.25%  "
